# Updates the cryptocurrency price list with freshly scraped values:
# per-row Price (column D) and 1h Volume-change % (column E) updates,
# plus two rows where the underlying coin changed (row 22/23 swap
# Litecoin<->PancakeSwap with refreshed figures, and row 51 TheGraph ->
# Cronos). Column D is written via a temporary Text number format so
# Excel doesn't reinterpret price strings like "300.90" or "0.0937" as
# numbers (which would drop the trailing zero / use scientific
# notation); the style is reset to Normal afterwards so no stray
# cell-format attribute is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.951.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.250.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "300.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.77%  "
$ws.Range("E7").Value = "  -2.47%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.899"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.589.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.263.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.863.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.81%  "
$ws.Range("E21").Value = "  -1.98%  "
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +26.79%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("E24").Value = "  -3.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "232.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.81%  "
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("E31").Value = "  -4.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "175.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +15.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0376"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.242"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.03%  "
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.18%  "
$ws.Range("E48").Value = "  +5.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "108.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0991"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.27%  "
